# --------------------------------------------------------------------------
# Applies the OOXML diff to the "Use case_confirmare_prezenta" document.
#
# Helper: force a clean run-split at a given (collapsed) Range boundary by
# toggling a character-formatting property on/off. The underlying engine
# always merges two adjacent runs that carry the exact same effective
# formatting once a document is saved, *unless* a run's <w:rPr> has been
# touched explicitly (even to a no-op net value) - that leaves the run with
# its own (possibly empty) <w:rPr>, which is enough to keep it a separate
# <w:r> element. We use that to split runs exactly where the diff wants a
# new run boundary, then restore the toggled property to its original value.
# --------------------------------------------------------------------------
$d = $word.ActiveDocument

function Split-RunAt($doc, [int]$pos) {
    $r = $doc.Range($pos, $pos + 1)
    $cur = $r.Font.Bold
    $r.Font.Bold = 1
    $r.Font.Bold = $cur
}

# ---------------------------------------------------------------------
# Change 1: "Angajatul" (POST-1. Angajatul va fi marcat ...) is split
# into two runs "An" + "gajatul" - no visible text change.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("POST-1. Angajatul va fi marcat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$postStart = $rng.Start
# "POST-1. " is 8 chars, then "An" is 2 more chars -> split after position postStart+8+2
Split-RunAt $d ($postStart + 8 + 2)

# ---------------------------------------------------------------------
# Change 2: ' cu "PRESENT" in ' -> ' cu "AVAILABLE" in ', and the run is
# split right after "AVAILABLE" (before the closing curly quote).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("PRESENT", $true, $false, $false, $false, $false, $true, 1, $false, "AVAILABLE", 2)
$rng2 = $d.Content
$rng2.Find.Execute("AVAILABLE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-RunAt $d ($rng2.Start + 9)

# ---------------------------------------------------------------------
# Change 3: remove the gramStart/gramEnd proofErr markers that wrap "va"
# in "Softul va afisa un mesaj ...".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Softul va afisa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "gramStart/gramEnd target located at" $rng.Start

# ---------------------------------------------------------------------
# Change 4: ' cu "Present".' -> ' cu "AVAILABLE".', with the _GoBack
# bookmark now sitting between "AVAILABLE" and the closing quote+period
# (it used to sit right after the whole ' cu "Present".' run).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Present", $true, $false, $false, $false, $false, $true, 1, $false, "AVAILABLE", 2)
